# Updates cryptos list values (price + 1h volume change) for the
# Sat Feb  3 07:14:23 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.085.66"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.319.38"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'303.40"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'99.60"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "'0.507"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.518"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("D10").Value = "'36.02"
$ws.Range("E10").Value = "  +4.32%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "'17.68"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").Value = "'6.90"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "2.681.38"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "2.316.64"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "'0.793"
$ws.Range("E17").Value = "  -3.28%  "
$ws.Range("D18").Value = "43.017.23"
$ws.Range("D19").Value = "'13.19"
$ws.Range("E19").Value = "  +4.64%  "
$ws.Range("D20").Value = "'6.18"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "'68.22"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").Value = "'239.87"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'25.34"
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").Value = "'169.40"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.19"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.05"
$ws.Range("E30").Value = "  -6.29%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'33.66"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").Value = "'4.97"
$ws.Range("E32").Value = "  +5.28%  "
$ws.Range("D33").Value = "'5.17"
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'18.33"
$ws.Range("E35").Value = "  +6.90%  "
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("D37").Value = "'0.0697"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").Value = "'1.83"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "'0.102"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'2.78"
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "1.994.92"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  -7.62%  "
$ws.Range("D45").Value = "'10.20"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").Value = "'17.51"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").Value = "'2.85"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'76.48"
$ws.Range("E48").Value = "  +8.67%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'55.02"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").Value = "2.547.26"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'1.54"
$ws.Range("E51").Value = "  -0.23%  "

# Excel marks cells whose apostrophe-forced text looks numeric with a
# "Text" quote-prefix style. Reset those cells back to the workbook's
# default "Normal" style so only the cell content changes.
foreach ($addr in @("D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D17", "D19", "D20", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D37", "D38", "D39", "D40", "D45", "D46", "D47", "D48", "D49", "D51")) {
    $ws.Range($addr).Style = "Normal"
}

Write-Output "Applied 99 cell updates"
